$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3820
$ws.Range("J29").Value = 4275
$ws.Range("L29").Value = 12825
$ws.Range("N29").Value = -13387

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3555.5217
$ws.Range("I40").Value = 1912.6428
$ws.Range("K40").Value = 1912.6428
$ws.Range("M40").Value = -1737.6428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4567.1
$ws.Range("I64").Value = 4500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4252

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4567.1
$ws.Range("I67").Value = 4500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3642

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4213.391
$ws.Range("I70").Value = 1329.1428
$ws.Range("J70").Value = 8700
$ws.Range("K70").Value = 3987.4284
$ws.Range("L70").Value = 26100
$ws.Range("M70").Value = -3717.4284
$ws.Range("N70").Value = -26640

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 4213.391
$ws.Range("I73").Value = 1329.1428
$ws.Range("J73").Value = 8700
$ws.Range("K73").Value = 3987.4284
$ws.Range("L73").Value = 26100
$ws.Range("M73").Value = -3051.4284
$ws.Range("N73").Value = -27972

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1833.3334
$ws.Range("I107").Value = 1833.3334
$ws.Range("K107").Value = 1833.3334
$ws.Range("M107").Value = 86.66660000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4131.1665
$ws.Range("I113").Value = 4435.769
$ws.Range("J113").Value = 3339.2
$ws.Range("K113").Value = 4435.769
$ws.Range("L113").Value = 3339.2
$ws.Range("M113").Value = -1181.769
$ws.Range("N113").Value = -9847.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1847.85
$ws.Range("J138").Value = 2499.2307
$ws.Range("L138").Value = 7497.6921
$ws.Range("N138").Value = -17777.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 31251938
$ws.Range("I74").Value = 38462916
$ws.Range("J74").Value = 4365.8335
$ws.Range("K74").Value = 38462916
$ws.Range("L74").Value = 4365.8335
$ws.Range("M74").Value = -38462042
$ws.Range("N74").Value = -6113.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 31251938
$ws.Range("I77").Value = 38462916
$ws.Range("J77").Value = 4365.8335
$ws.Range("K77").Value = 192314580
$ws.Range("L77").Value = 21829.1675
$ws.Range("M77").Value = -192310212
$ws.Range("N77").Value = -30565.1675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 889.5
$ws.Range("I97").Value = 1037.0625
$ws.Range("K97").Value = 1037.0625
$ws.Range("M97").Value = -541.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1507.55
$ws.Range("I102").Value = 1261.8125
$ws.Range("K102").Value = 1261.8125
$ws.Range("M102").Value = 360.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 143579.28
$ws.Range("I110").Value = 143579.28
$ws.Range("K110").Value = 143579.28
$ws.Range("M110").Value = -141534.28

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 78999
$ws.Range("J6").Value = 78999
$ws.Range("L6").Value = 78999
$ws.Range("N6").Value = -79225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1915.6
$ws.Range("I20").Value = 1745
$ws.Range("J20").Value = 2313.6667
$ws.Range("K20").Value = 1745
$ws.Range("L20").Value = 2313.6667
$ws.Range("M20").Value = -1498
$ws.Range("N20").Value = -2807.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 667.3
$ws.Range("J80").Value = 758
$ws.Range("L80").Value = 758
$ws.Range("N80").Value = -2754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 667.3
$ws.Range("J83").Value = 758
$ws.Range("L83").Value = 3790
$ws.Range("N83").Value = -13774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 112499
$ws.Range("J114").Value = 112499
$ws.Range("L114").Value = 112499
$ws.Range("N114").Value = -121177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27785476
$ws.Range("I58").Value = 71443970
$ws.Range("K58").Value = 71443970
$ws.Range("M58").Value = -71443767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1045.9375
$ws.Range("I105").Value = 1060.7858
$ws.Range("J105").Value = 942
$ws.Range("K105").Value = 1060.7858
$ws.Range("L105").Value = 942
$ws.Range("M105").Value = 686.2141999999999
$ws.Range("N105").Value = -4436

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 27780274
$ws.Range("I134").Value = 35715852
$ws.Range("K134").Value = 107147556
$ws.Range("M134").Value = -107145021

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 27785476
$ws.Range("I136").Value = 71443970
$ws.Range("K136").Value = 214331910
$ws.Range("M136").Value = -214329360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7846.5
$ws.Range("J80").Value = 7796.5
$ws.Range("L80").Value = 23389.5
$ws.Range("N80").Value = -25261.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 7846.5
$ws.Range("J83").Value = 7796.5
$ws.Range("L83").Value = 70168.5
$ws.Range("N83").Value = -79528.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 731.0526
$ws.Range("I107").Value = 187.625
$ws.Range("J107").Value = 1126.2727
$ws.Range("K107").Value = 562.875
$ws.Range("L107").Value = 3378.8181
$ws.Range("M107").Value = 1357.125
$ws.Range("N107").Value = -7218.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 74999
$ws.Range("J110").Value = 74999
$ws.Range("L110").Value = 74999
$ws.Range("N110").Value = -83179

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 50962.383
$ws.Range("I113").Value = 59083.223
$ws.Range("J113").Value = 2237.3333
$ws.Range("K113").Value = 59083.223
$ws.Range("L113").Value = 2237.3333
$ws.Range("M113").Value = -56913.223
$ws.Range("N113").Value = -6577.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3334683
$ws.Range("I46").Value = 1619.6
$ws.Range("K46").Value = 1619.6
$ws.Range("M46").Value = -1431.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3370
$ws.Range("I81").Value = 2545.4443
$ws.Range("K81").Value = 5090.8886
$ws.Range("M81").Value = -4029.8886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3370
$ws.Range("I84").Value = 2545.4443
$ws.Range("K84").Value = 25454.443
$ws.Range("M84").Value = -20150.443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 115275
$ws.Range("J92").Value = 115275
$ws.Range("L92").Value = 115275
$ws.Range("N92").Value = -120267

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 46266.668
$ws.Range("J95").Value = 46266.668
$ws.Range("L95").Value = 46266.668
$ws.Range("N95").Value = -51758.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 984.36664
$ws.Range("I113").Value = 965.2083
$ws.Range("J113").Value = 1061
$ws.Range("K113").Value = 2895.6249
$ws.Range("L113").Value = 3183
$ws.Range("M113").Value = -725.6248999999998
$ws.Range("N113").Value = -7523

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13515023
$ws.Range("I136").Value = 16668149
$ws.Range("J136").Value = 1625
$ws.Range("K136").Value = 50004447
$ws.Range("L136").Value = 4875
$ws.Range("M136").Value = -50001897
$ws.Range("N136").Value = -9975
